# Updated queries for C3DC first half testcases.
# Replace the stale `std.id` / `prt.id` join keys with the new
# `study_id` / `participant_id` join keys across every SQL query
# stored on Sheet1 (cells B2:B7 and C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldJoin = "LEFT JOIN `n    df_participant prt ON std.id = prt.`"study.id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.id = dgn.`"participant.id`"`nLEFT JOIN `n    df_treatments trt ON prt.id = trt.`"participant.id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.id = trr.`"participant.id`"`nLEFT JOIN `n    df_survival srv ON prt.id = srv.`"participant.id`"`nLEFT JOIN `n    df_reference_files rfs ON std.id = rfs.`"study.id`""

$newJoin = "LEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`""

$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")
foreach ($cellRef in $cells) {
    $rng = $ws.Range($cellRef)
    $current = $rng.Value()
    $rng.Value = $current.Replace($oldJoin, $newJoin)
}

# The author scrolled the sheet down and left the cursor on C7 before saving.
$ws.Range("C7").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

# Widen column C (StatQuery) to fit the longer join text; this also drops
# the old "bestFit" autofit flag since the width is now explicit.
$ws.Columns.Item(3).ColumnWidth = 72
